$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1454.2273
$ws.Range("I33").Value = 1228.9474
$ws.Range("K33").Value = 1228.9474
$ws.Range("M33").Value = -999.9474
$ws.Range("H117").Value = 132969.5
$ws.Range("J117").Value = 132969.5
$ws.Range("L117").Value = 132969.5
$ws.Range("N117").Value = -142147.5
$ws.Range("H127").Value = 2374.5
$ws.Range("J127").Value = 2749.5
$ws.Range("L127").Value = 8248.5
$ws.Range("N127").Value = -18168.5
$ws.Range("H129").Value = 1566.1666
$ws.Range("I129").Value = 560
$ws.Range("J129").Value = 2572.3333
$ws.Range("K129").Value = 1680
$ws.Range("L129").Value = 7716.999899999999
$ws.Range("M129").Value = 3320
$ws.Range("N129").Value = -17716.9999
$ws.Range("H138").Value = 3046.6829
$ws.Range("J138").Value = 3208.6316
$ws.Range("L138").Value = 9625.8948
$ws.Range("N138").Value = -19905.8948
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 29776.666
$ws.Range("I43").Value = 27998.5
$ws.Range("K43").Value = 27998.5
$ws.Range("M43").Value = -27685.5
$ws.Range("H45").Value = 1617.2858
$ws.Range("I45").Value = 1427.5555
$ws.Range("K45").Value = 1427.5555
$ws.Range("M45").Value = -1050.5555
$ws.Range("H74").Value = 10007756
$ws.Range("I74").Value = 16667398
$ws.Range("J74").Value = 18293.7
$ws.Range("K74").Value = 16667398
$ws.Range("L74").Value = 18293.7
$ws.Range("M74").Value = -16666524
$ws.Range("N74").Value = -20041.7
$ws.Range("H77").Value = 10007756
$ws.Range("I77").Value = 16667398
$ws.Range("J77").Value = 18293.7
$ws.Range("K77").Value = 83336990
$ws.Range("L77").Value = 91468.5
$ws.Range("M77").Value = -83332622
$ws.Range("N77").Value = -100204.5
$ws.Range("H122").Value = 1304.3077
$ws.Range("I122").Value = 1168.1
$ws.Range("J122").Value = 1758.3334
$ws.Range("K122").Value = 3504.3
$ws.Range("L122").Value = 5275.0002
$ws.Range("M122").Value = -1054.3
$ws.Range("N122").Value = -10175.0002
$ws.Range("H132").Value = 3707.4106
$ws.Range("I132").Value = 2436.925
$ws.Range("J132").Value = 6883.625
$ws.Range("K132").Value = 7310.775000000001
$ws.Range("L132").Value = 20650.875
$ws.Range("M132").Value = -4780.775000000001
$ws.Range("N132").Value = -25710.875
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1081.1562
$ws.Range("I94").Value = 1089.5807
$ws.Range("K94").Value = 1089.5807
$ws.Range("M94").Value = -638.5807
$ws.Range("H105").Value = 2380.5
$ws.Range("I105").Value = 1561
$ws.Range("K105").Value = 1561
$ws.Range("M105").Value = 186
$ws.Range("H107").Value = 2336.3333
$ws.Range("I107").Value = 2336.3333
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 2336.3333
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -416.3332999999998
$ws.Range("N107").ClearContents()
$ws.Range("H134").Value = 40711.816
$ws.Range("I134").Value = 2060.2104
$ws.Range("J134").Value = 132509.38
$ws.Range("K134").Value = 6180.6312
$ws.Range("L134").Value = 397528.14
$ws.Range("M134").Value = -3645.6312
$ws.Range("N134").Value = -402598.14
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H64").Value = 109000
$ws.Range("J64").Value = 109000
$ws.Range("L64").Value = 109000
$ws.Range("N64").Value = -109496
$ws.Range("H67").Value = 109000
$ws.Range("J67").Value = 109000
$ws.Range("L67").Value = 109000
$ws.Range("N67").Value = -110716
$ws.Range("H70").Value = 47500
$ws.Range("I70").Value = 40000
$ws.Range("K70").Value = 40000
$ws.Range("M70").Value = -39685
$ws.Range("H73").Value = 47500
$ws.Range("I73").Value = 40000
$ws.Range("K73").Value = 40000
$ws.Range("M73").Value = -38908
$ws.Range("H88").Value = 39447.668
$ws.Range("J88").Value = 39447.668
$ws.Range("L88").Value = 39447.668
$ws.Range("N88").Value = -40259.668
$ws.Range("H91").Value = 39447.668
$ws.Range("J91").Value = 39447.668
$ws.Range("L91").Value = 39447.668
$ws.Range("N91").Value = -42255.668
$ws.Range("H94").Value = 8045.6665
$ws.Range("I94").Value = 7867
$ws.Range("J94").Value = 8135
$ws.Range("K94").Value = 7867
$ws.Range("L94").Value = 8135
$ws.Range("M94").Value = -7416
$ws.Range("N94").Value = -9037
$ws.Range("H105").Value = 2136.75
$ws.Range("I105").Value = 2009.6666
$ws.Range("J105").Value = 2213
$ws.Range("K105").Value = 2009.6666
$ws.Range("L105").Value = 2213
$ws.Range("M105").Value = -262.6666
$ws.Range("N105").Value = -5707
$ws.Range("H132").Value = 1994.8889
$ws.Range("I132").Value = 1994.8889
$ws.Range("K132").Value = 5984.6667
$ws.Range("M132").Value = -3454.6667
$ws.Range("H134").Value = 1119249.5
$ws.Range("I134").Value = 2501558.8
$ws.Range("K134").Value = 7504676.399999999
$ws.Range("M134").Value = -7502141.399999999
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 5154.6
$ws.Range("J80").Value = 4943.625
$ws.Range("L80").Value = 14830.875
$ws.Range("N80").Value = -16702.875
$ws.Range("H83").Value = 5154.6
$ws.Range("J83").Value = 4943.625
$ws.Range("L83").Value = 44492.625
$ws.Range("N83").Value = -53852.625
$ws.Range("H129").Value = 1127.875
$ws.Range("I129").Value = 812.3333
$ws.Range("J129").Value = 2074.5
$ws.Range("K129").Value = 2436.9999
$ws.Range("L129").Value = 6223.5
$ws.Range("M129").Value = 2563.0001
$ws.Range("N129").Value = -16223.5
$ws.Range("H137").Value = 5462.643
$ws.Range("J137").Value = 5247.5
$ws.Range("L137").Value = 15742.5
$ws.Range("N137").Value = -25942.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3615.7083
$ws.Range("I102").Value = 2671.7727
$ws.Range("K102").Value = 2671.7727
$ws.Range("M102").Value = -1049.7727
$ws.Range("H132").Value = 40002590
$ws.Range("I132").Value = 52634108
$ws.Range("K132").Value = 157902324
$ws.Range("M132").Value = -157899794
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 12567080
$ws.Range("I7").Value = 18185434
$ws.Range("K7").Value = 18185434
$ws.Range("M7").Value = -18185322
$ws.Range("H40").Value = 2923.6365
$ws.Range("I40").Value = 2020.1875
$ws.Range("K40").Value = 2020.1875
$ws.Range("M40").Value = -1884.1875
$ws.Range("H94").Value = 55000
$ws.Range("J94").Value = 55000
$ws.Range("L94").Value = 55000
$ws.Range("N94").Value = -56352
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()
$ws.Range("H122").Value = 4408.7144
$ws.Range("I122").Value = 3721.276
$ws.Range("J122").Value = 7731.3335
$ws.Range("K122").Value = 11163.828
$ws.Range("L122").Value = 23194.0005
$ws.Range("M122").Value = -8713.828
$ws.Range("N122").Value = -28094.0005
$ws.Range("H126").Value = 12567080
$ws.Range("I126").Value = 18185434
$ws.Range("K126").Value = 54556302
$ws.Range("M126").Value = -54553832
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H19").Value = 3000
$ws.Range("J19").Value = 3000
$ws.Range("L19").Value = 3000
$ws.Range("N19").Value = -3348
$ws.Range("H113").Value = 1143.4667
$ws.Range("I113").Value = 1182.75
$ws.Range("K113").Value = 3548.25
$ws.Range("M113").Value = -1378.25
$ws.Range("H132").Value = 1448.125
$ws.Range("I132").Value = 1467.6522
$ws.Range("J132").Value = 999
$ws.Range("K132").Value = 4402.9566
$ws.Range("L132").Value = 2997
$ws.Range("M132").Value = -1872.9566
$ws.Range("N132").Value = -8057
